$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for data rows 2-221 changes from 2023-10-03 (45202)
# to 2023-10-04 (45203) for every row. Use the raw Excel date serial number
# so no time-of-day component gets attached to the value.
$ws.Range("C2:C221").Value = 45203
